$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 16 - fill in with the same PUBLONS010 data already used elsewhere in the sheet
$ws.Range("A16").Value = "PUBLONS010"
$ws.Range("B16").Value = "OPQA-5784&&OPQA-5785"
$ws.Range("C16").Value = "Verify that `"Your email address is already registered. Please sign in.`" error message whenever try to create publons user using existing account.&&Verify that email address prepopulated in sign in page whenever try to register user using existing user"
$ws.Range("D16").Value = "Y"

# Row 17 - new script PUBLONS011 / OPQA-5986
$ws.Range("A17").Value = "PUBLONS011"
$ws.Range("B17").Value = "OPQA-5986"
$ws.Range("C17").Value = "Verify that Customer care component links for all applications`nCMTY - 'community.info@clarivate.com' linked to community.info@clarivate.com email`nCADP - 'Customer Support' linked to https://support.clarivate.com/s/`nDRA - 'Drug Research Advisor Customer Care' linked to https://support.clarivate.com/LifeSciences/`nCMC - 'Cortellis CMC Intelligence Customer Care' linked to https://support.clarivate.com/LifeSciences/`nINTEGRITY - 'Integrity Customer Care' linked to https://support.clarivate.com/LifeSciences/`nENDNOTE - 'EndNote Customer Care' linked to http://endnote.com/support`nIPA - 'IPA.support@thomsonreuters.com' linked to IPA.support@thomsonreuters.com`nPUBLONS - 'Customer Care' linked to info@publons.com`nPSA - 'community.info@clarivate.com' linked to community.info@clarivate.com`nWAT - 'sarlabs.info@clarivate.com' linked to sarlabs.info@clarivate.com`n"
$ws.Range("D17").Value = "Y"

# Row 18 - new script PUBLONS012 / OPQA-5859&&OPQA-5860
$ws.Range("A18").Value = "PUBLONS012"
$ws.Range("B18").Value = "OPQA-5859&&OPQA-5860"
$ws.Range("C18").Value = "Verify that TERMS OF USE and PRIVACY STATEMENT links are working correctly in publons landing page && Verify that TERMS OF USE and PRIVACY STATEMENT links are working correctly in publons Registration page"
$ws.Range("D18").Value = "Y"

# Row 19 - new script (JIRA + description entered before the TCID)
$ws.Range("B19").Value = "OPQA-5984"
$ws.Range("C19").Value = "Verify that Email should be transferred when switching from sign in to register for CMTY,PUBLONS,ENDNOTE"
$ws.Range("D19").Value = "Y"
$ws.Range("A19").Value = "PUBLONS017"

# Row heights follow the wrapped content in column C
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 195
$ws.Rows.Item(18).RowHeight = 30

# Scroll the view down and move the selection onto the newly-added row
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A19").Select()
